# Auto-generated Excel COM-interop script applying the Titan_Profits workbook update.
# For each affected row (columns H:N — price/profit calc columns), sets the new
# values and clears any cell that the diff removes entirely.

$wb = $excel.ActiveWorkbook

# ALC!row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1947
$ws.Range("I51").Value = 1438.1666
$ws.Range("J51").Value = 5000
$ws.Range("K51").Value = 1438.1666
$ws.Range("L51").Value = 5000
$ws.Range("M51").Value = -954.1666
$ws.Range("N51").Value = -5968

# ALC!row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1600.5264
$ws.Range("I70").Value = 1323.5
$ws.Range("J70").Value = 1674.4
$ws.Range("K70").Value = 3970.5
$ws.Range("L70").Value = 5023.200000000001
$ws.Range("M70").Value = -3700.5
$ws.Range("N70").Value = -5563.200000000001

# ALC!row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1600.5264
$ws.Range("I73").Value = 1323.5
$ws.Range("J73").Value = 1674.4
$ws.Range("K73").Value = 3970.5
$ws.Range("L73").Value = 5023.200000000001
$ws.Range("M73").Value = -3034.5
$ws.Range("N73").Value = -6895.200000000001

# ALC!row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1200.4286
$ws.Range("I80").Value = 4
$ws.Range("J80").Value = 1399.8334
$ws.Range("K80").Value = 12
$ws.Range("L80").Value = 4199.5002
$ws.Range("M80").Value = 986
$ws.Range("N80").Value = -6195.5002

# ALC!row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1200.4286
$ws.Range("I83").Value = 4
$ws.Range("J83").Value = 1399.8334
$ws.Range("K83").Value = 36
$ws.Range("L83").Value = 12598.5006
$ws.Range("M83").Value = 4956
$ws.Range("N83").Value = -22582.5006

# ALC!row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 700127.4
$ws.Range("I98").Value = 931711
$ws.Range("J98").Value = 5376.5
$ws.Range("K98").Value = 931711
$ws.Range("L98").Value = 5376.5
$ws.Range("M98").Value = -930213
$ws.Range("N98").Value = -8372.5

# ALC!row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 700127.4
$ws.Range("I122").Value = 931711
$ws.Range("J122").Value = 5376.5
$ws.Range("K122").Value = 2795133
$ws.Range("L122").Value = 16129.5
$ws.Range("M122").Value = -2792683
$ws.Range("N122").Value = -21029.5

# ALC!row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2187.25
$ws.Range("J129").Value = 7000
$ws.Range("L129").Value = 21000
$ws.Range("N129").Value = -31000

# ALC!row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 415622.38
$ws.Range("I132").Value = 486406.6
$ws.Range("J132").Value = 61701.2
$ws.Range("K132").Value = 1459219.8
$ws.Range("L132").Value = 185103.6
$ws.Range("M132").Value = -1456689.8
$ws.Range("N132").Value = -190163.6

# ARM!row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8850.733
$ws.Range("I2").Value = 11635.546
$ws.Range("K2").Value = 11635.546
$ws.Range("M2").Value = -11522.546

# ARM!row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2664.182
$ws.Range("I61").Value = 1868.1875
$ws.Range("J61").Value = 4786.8335
$ws.Range("K61").Value = 1868.1875
$ws.Range("L61").Value = 4786.8335
$ws.Range("M61").Value = -1656.1875
$ws.Range("N61").Value = -5210.8335

# ARM!row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 8850.733
$ws.Range("I116").Value = 11635.546
$ws.Range("K116").Value = 11635.546
$ws.Range("M116").Value = -9341.546

# ARM!row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2664.182
$ws.Range("I136").Value = 1868.1875
$ws.Range("J136").Value = 4786.8335
$ws.Range("K136").Value = 5604.5625
$ws.Range("L136").Value = 14360.5005
$ws.Range("M136").Value = -3054.5625
$ws.Range("N136").Value = -19460.5005

# ARM!row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 65395.4
$ws.Range("J139").Value = 65395.4
$ws.Range("L139").Value = 65395.4
$ws.Range("N139").Value = -75675.39999999999

# ARM!row 140
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

# ARM!row 141
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 81045.8
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 81045.8
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 81045.8
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -91405.8

# BSM!row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8850.733
$ws.Range("I3").Value = 11635.546
$ws.Range("K3").Value = 11635.546
$ws.Range("M3").Value = -11521.546

# BSM!row 74
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 41000
$ws.Range("J74").Value = 41000
$ws.Range("L74").Value = 41000
$ws.Range("N74").Value = -42872

# BSM!row 77
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H77").Value = 41000
$ws.Range("J77").Value = 41000
$ws.Range("L77").Value = 123000
$ws.Range("N77").Value = -132360

# BSM!row 106
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 28668.2
$ws.Range("J106").Value = 28668.2
$ws.Range("L106").Value = 28668.2
$ws.Range("N106").Value = -31192.2

# CRP!row 13
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 522500
$ws.Range("J13").Value = 522500
$ws.Range("L13").Value = 522500
$ws.Range("N13").Value = -522778

# CRP!row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2074.5
$ws.Range("I86").Value = 1700
$ws.Range("K86").Value = 1700
$ws.Range("M86").Value = -577

# CRP!row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 2074.5
$ws.Range("I89").Value = 1700
$ws.Range("K89").Value = 8500
$ws.Range("M89").Value = -2884

# CRP!row 103
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 17069.4
$ws.Range("I103").Value = 12636.75
$ws.Range("J103").Value = 34800
$ws.Range("K103").Value = 12636.75
$ws.Range("L103").Value = 34800
$ws.Range("M103").Value = -11464.75
$ws.Range("N103").Value = -37144

# CRP!row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 809
$ws.Range("I122").Value = 809
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2427
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 23
$ws.Range("N122").ClearContents()

# CUL!row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1235.9615
$ws.Range("I5").Value = 876.3333
$ws.Range("J5").Value = 1343.85
$ws.Range("K5").Value = 2628.9999
$ws.Range("L5").Value = 4031.55
$ws.Range("M5").Value = -2516.9999
$ws.Range("N5").Value = -4255.549999999999

# CUL!row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 627.34784
$ws.Range("J122").Value = 1071.4
$ws.Range("L122").Value = 9642.6
$ws.Range("N122").Value = -14542.6

# CUL!row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1089.5
$ws.Range("I131").Value = 468.8889
$ws.Range("J131").Value = 1150.8792
$ws.Range("K131").Value = 1406.6667
$ws.Range("L131").Value = 3452.6376
$ws.Range("M131").Value = 3633.3333
$ws.Range("N131").Value = -13532.6376

# CUL!row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1235.9615
$ws.Range("I135").Value = 876.3333
$ws.Range("J135").Value = 1343.85
$ws.Range("K135").Value = 7886.9997
$ws.Range("L135").Value = 12094.65
$ws.Range("M135").Value = -5351.9997
$ws.Range("N135").Value = -17164.65

# GSM!row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1991.8
$ws.Range("J102").Value = 2029.5714
$ws.Range("L102").Value = 2029.5714
$ws.Range("N102").Value = -5273.5714

# LTW!row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 775.94116
$ws.Range("I22").Value = 782.4167
$ws.Range("K22").Value = 782.4167
$ws.Range("M22").Value = -487.4167

# LTW!row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 775.94116
$ws.Range("I27").Value = 782.4167
$ws.Range("K27").Value = 782.4167
$ws.Range("M27").Value = -675.4167

# LTW!row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1054.7142
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
